# Natmi following Dr Hou advice
#
# The NATMI ligand-receptor edge table for Efna5 -> Ephb1 was recomputed
# (rows 2-3 values changed) and two additional rows were added for the
# newly-included "sCs" sending/target cluster, giving four
# sending-cluster/target-cluster combinations in total:
#   FAPs -> ECs, FAPs -> sCs, sCs -> ECs, sCs -> sCs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns A-D: Sending cluster, Ligand symbol, Receptor symbol, Target cluster
# Columns E-T: the 16 NATMI numeric statistics columns

$data = @(
    @{ Row = 2;  A = "FAPs"; B = "Efna5"; C = "Ephb1"; D = "ECs";
       Vals = @(3, 1, 2.030023666666667, 6.090071, 0.8776223887075381, 0.8776223887075382,
                3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272,
                4.490141680844445, 40.4112751276, 0.4222432562762945, 0.4222432562762946) },
    @{ Row = 3;  A = "FAPs"; B = "Efna5"; C = "Ephb1"; D = "sCs";
       Vals = @(3, 1, 2.030023666666667, 6.090071, 0.8776223887075381, 0.8776223887075382,
                3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727,
                4.84250912885711, 43.58258215971399, 0.4553791324312435, 0.4553791324312436) },
    @{ Row = 4;  A = "sCs";  B = "Efna5"; C = "Ephb1"; D = "ECs";
       Vals = @(2, 0.6666666666666666, 0.283071, 0.849213, 0.1223776112924619, 0.1223776112924619,
                3, 1, 2.211866666666667, 6.6356, 0.4811217919111272, 0.4811217919111272,
                0.6261153092000001, 5.6350377828, 0.05887853563483265, 0.05887853563483265) },
    @{ Row = 5;  A = "sCs";  B = "Efna5"; C = "Ephb1"; D = "sCs";
       Vals = @(2, 0.6666666666666666, 0.283071, 0.849213, 0.1223776112924619, 0.1223776112924619,
                3, 1, 2.385444666666666, 7.156333999999999, 0.5188782080888727, 0.5188782080888727,
                0.675250207238, 6.077251865141999, 0.06349907565762922, 0.06349907565762922) }
)

$numCols = @("E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T")

foreach ($entry in $data) {
    $r = $entry.Row
    $ws.Range("A$r").Value = $entry.A
    $ws.Range("B$r").Value = $entry.B
    $ws.Range("C$r").Value = $entry.C
    $ws.Range("D$r").Value = $entry.D

    for ($i = 0; $i -lt $numCols.Count; $i++) {
        $ws.Range("$($numCols[$i])$r").Value = $entry.Vals[$i]
    }
}
